$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'62.779.76"
$ws.Range("E2").Value = "'  +9.80%  "

# Row 3
$ws.Range("D3").Value = "'3.463.06"
$ws.Range("E3").Value = "'  +6.88%  "

# Row 4
$ws.Range("D4").Value = "'1.01"
$ws.Range("E4").Value = "'  +0.50%  "

# Row 5
$ws.Range("D5").Value = "'421.71"
$ws.Range("E5").Value = "'  +6.68%  "

# Row 6
$ws.Range("D6").Value = "'114.55"
$ws.Range("E6").Value = "'  +6.69%  "

# Row 7
$ws.Range("D7").Value = "'0.595"
$ws.Range("E7").Value = "'  +4.50%  "

# Row 8
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "'  +0.27%  "

# Row 9
$ws.Range("D9").Value = "'0.651"
$ws.Range("E9").Value = "'  +5.70%  "

# Row 10
$ws.Range("D10").Value = "'0.114"
$ws.Range("E10").Value = "'  +18.23%  "

# Row 11
$ws.Range("D11").Value = "'40.40"
$ws.Range("E11").Value = "'  +3.78%  "

# Row 12
$ws.Range("E12").Value = "'  +1.30%  "

# Row 13
$ws.Range("D13").Value = "'4.020.19"
$ws.Range("E13").Value = "'  +7.07%  "

# Row 14
$ws.Range("D14").Value = "'8.51"
$ws.Range("E14").Value = "'  +5.01%  "

# Row 15
$ws.Range("D15").Value = "'19.98"
$ws.Range("E15").Value = "'  +5.64%  "

# Row 16
$ws.Range("D16").Value = "'3.517.74"
$ws.Range("E16").Value = "'  +8.85%  "

# Row 17
$ws.Range("B17").Value = "'WrappedBTC"
$ws.Range("C17").Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "'63.040.73"
$ws.Range("E17").Value = "'  +10.65%  "

# Row 18
$ws.Range("B18").Value = "'Polygon"
$ws.Range("C18").Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D18").Value = "'1.05"
$ws.Range("E18").Value = "'  +1.67%  "

# Row 19
$ws.Range("D19").Value = "'10.94"
$ws.Range("E19").Value = "'  -0.64%  "

# Row 20
$ws.Range("D20").Value = "'0.0000118"
$ws.Range("E20").Value = "'  +10.91%  "

# Row 21
$ws.Range("D21").Value = "'3.40"
$ws.Range("E21").Value = "'  +2.30%  "

# Row 22
$ws.Range("D22").Value = "'13.13"
$ws.Range("E22").Value = "'  +1.16%  "

# Row 23
$ws.Range("D23").Value = "'305.72"
$ws.Range("E23").Value = "'  +2.50%  "

# Row 24
$ws.Range("D24").Value = "'76.43"
$ws.Range("E24").Value = "'  +3.23%  "

# Row 25
$ws.Range("D25").Value = "'3.31"
$ws.Range("E25").Value = "'  +5.32%  "

# Row 26
$ws.Range("D26").Value = "'30.32"
$ws.Range("E26").Value = "'  +8.74%  "

# Row 27
$ws.Range("E27").Value = "'  +2.86%  "

# Row 28
$ws.Range("D28").Value = "'7.92"
$ws.Range("E28").Value = "'  +3.06%  "

# Row 29
$ws.Range("B29").Value = "'Kaspa"
$ws.Range("C29").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "'0.177"
$ws.Range("E29").Value = "'  +5.18%  "

# Row 30
$ws.Range("B30").Value = "'RenderToken"
$ws.Range("C30").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").Value = "'7.53"
$ws.Range("E30").Value = "'  +3.96%  "

# Row 31
$ws.Range("D31").Value = "'0.114"
$ws.Range("E31").Value = "'  +5.13%  "

# Row 32
$ws.Range("D32").Value = "'11.52"
$ws.Range("E32").Value = "'  +4.80%  "

# Row 33
$ws.Range("B33").Value = "'Dai"
$ws.Range("C33").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "'  -0.05%  "

# Row 34
$ws.Range("B34").Value = "'Toncoin"
$ws.Range("C34").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D34").Value = "'2.49"
$ws.Range("E34").Value = "'  +18.21%  "

# Row 35
$ws.Range("D35").Value = "'39.75"
$ws.Range("E35").Value = "'  +6.15%  "

# Row 36
$ws.Range("D36").Value = "'0.0509"
$ws.Range("E36").Value = "'  +5.24%  "

# Row 37
$ws.Range("D37").Value = "'52.00"
$ws.Range("E37").Value = "'  +0.60%  "

# Row 38
$ws.Range("D38").Value = "'3.13"
$ws.Range("E38").Value = "'  +3.42%  "

# Row 39
$ws.Range("B39").Value = "'FirstDigitalUSD"
$ws.Range("C39").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "'  +0.15%  "

# Row 40
$ws.Range("B40").Value = "'LidoDAOToken"
$ws.Range("C40").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "'3.42"
$ws.Range("E40").Value = "'  -2.89%  "

# Row 41
$ws.Range("B41").Value = "'Monero"
$ws.Range("C41").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "'138.92"
$ws.Range("E41").Value = "'  +3.13%  "

# Row 42
$ws.Range("B42").Value = "'Stellar"
$ws.Range("C42").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "'0.123"
$ws.Range("E42").Value = "'  +3.01%  "

# Row 43
$ws.Range("B43").Value = "'ARBITRUM"
$ws.Range("C43").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'1.94"
$ws.Range("E43").Value = "'  +2.49%  "

# Row 44
$ws.Range("B44").Value = "'TheGraph"
$ws.Range("C44").Value = "'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "'0.290"
$ws.Range("E44").Value = "'  +3.17%  "

# Row 45
$ws.Range("B45").Value = "'NEARProtocol"
$ws.Range("C45").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "'3.99"
$ws.Range("E45").Value = "'  +0.71%  "

# Row 46
$ws.Range("B46").Value = "'Celestia"
$ws.Range("C46").Value = "'https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D46").Value = "'16.86"
$ws.Range("E46").Value = "'  -0.19%  "

# Row 47
$ws.Range("B47").Value = "'WEMIXToken"
$ws.Range("C47").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "'2.29"
$ws.Range("E47").Value = "'  +9.68%  "

# Row 48
$ws.Range("B48").Value = "'EnergySwap"
$ws.Range("C48").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'22.62"
$ws.Range("E48").Value = "'  +3.46%  "

# Row 49
$ws.Range("D49").Value = "'2.202.29"
$ws.Range("E49").Value = "'  +2.54%  "

# Row 50
$ws.Range("B50").Value = "'ApeXProtocol"
$ws.Range("C50").Value = "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D50").Value = "'2.38"
$ws.Range("E50").Value = "'  +1.85%  "

# Row 51
$ws.Range("B51").Value = "'ThetaToken"
$ws.Range("C51").Value = "'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "'1.97"
$ws.Range("E51").Value = "'  -2.89%  "
